# feat: add 2022-Q1 data
#
# 1) Insert a new sheet "2022-Q1" (fund-level holdings) positioned right
#    before the "总计" (totals) sheet.
# 2) Prepend a new summary row to "总计" for 2022-Q1 and renumber the
#    zero-based index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: force a numeric-looking string ("013233", "39.48", ...) to be
# stored as genuine text instead of being auto-coerced to a number by
# the normal .Value setter. We do this by writing a text formula
# (="013233") and then collapsing it down to a plain value with a
# values-only paste, which keeps the existing (unstyled) cell format.
# ---------------------------------------------------------------------
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# =======================================================================
# Part 1: new "2022-Q1" worksheet
# =======================================================================

$wsTotal = $wb.Worksheets.Item("总计")

# Duplicate "总计" so the new sheet inherits identical sheetPr / page
# setup / formatting, then drop it in right before "总计" and rename.
$wsTotal.Copy($wsTotal)
$wsQ1 = $wb.Worksheets.Item(5)
$wsQ1.Name = "2022-Q1"

# Header row (B1:H1) - B1:D1 already carry the header style from the
# copied sheet; extend the same style across the new E1:H1 cells.
$wsQ1.Cells.Item(1,2).Value = "基金代码"
$wsQ1.Cells.Item(1,3).Value = "基金名称"
$wsQ1.Cells.Item(1,4).Value = "基金规模"
$wsQ1.Cells.Item(1,5).Value = "股票总仓位"
$wsQ1.Cells.Item(1,6).Value = "仓位占比"
$wsQ1.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ1.Cells.Item(1,8).Value = "仓位排名"

$wsQ1.Cells.Item(1,2).Copy()
$wsQ1.Range("E1:H1").PasteSpecial(-4122)

# Data rows: index, code, name, size, total position, position ratio,
# held value, position rank.
$rows = @(
    @(0, "013233", "华夏中证500指数智选增强A",   "39.48", "92.73", "1.39", "0.5488", 7),
    @(1, "005994", "国投瑞银中证500指数量化增强A", "11.53", "87.00", "1.09", "0.1257", 10),
    @(2, "013234", "华夏中证500指数智选增强C",   "4.28",  "92.73", "1.39", "0.0595", 7),
    @(3, "007089", "国投瑞银中证500指数量化增强C", "3.82",  "87.00", "1.09", "0.0416", 10),
    @(4, "005536", "渤海汇金量化成长混合",        "0.61",  "88.57", "0.75", "0.0046", 7)
)

# Copy the existing A-column index style onto row 6 (a brand new row
# beyond the 5 rows the source "总计" sheet had) before we start writing
# values into it.
$wsQ1.Cells.Item(2,1).Copy()
$wsQ1.Cells.Item(6,1).PasteSpecial(-4122)

foreach ($r in $rows) {
    $rowNum = [int]$r[0] + 2
    $wsQ1.Cells.Item($rowNum, 1).Value = [int]$r[0]
    Set-TextValue $wsQ1.Cells.Item($rowNum, 2) $r[1]
    $wsQ1.Cells.Item($rowNum, 3).Value = $r[2]
    Set-TextValue $wsQ1.Cells.Item($rowNum, 4) $r[3]
    Set-TextValue $wsQ1.Cells.Item($rowNum, 5) $r[4]
    Set-TextValue $wsQ1.Cells.Item($rowNum, 6) $r[5]
    Set-TextValue $wsQ1.Cells.Item($rowNum, 7) $r[6]
    $wsQ1.Cells.Item($rowNum, 8).Value = $r[7]
}

# =======================================================================
# Part 2: prepend the 2022-Q1 row to "总计"
# =======================================================================

# IMPORTANT: worksheet handles in this host track by *position*, not a
# stable object identity. $wsTotal was captured before the sheet-copy
# above shifted "总计" from slot 5 to slot 6, so it would now silently
# resolve to the new "2022-Q1" sheet instead. Re-resolve it by name.
$wsTotal = $wb.Worksheets.Item("总计")

# Apply the existing index-column style to the new bottom row (6) before
# shifting values down into it.
$wsTotal.Cells.Item(2,1).Copy()
$wsTotal.Cells.Item(6,1).PasteSpecial(-4122)

# Shift rows 2-5 down to 3-6 (bottom-up so we don't clobber a row before
# reading it). Use .Value2 to read back the literal stored value - plain
# .Value as an rvalue does not reliably round-trip through this host.
$wsTotal.Cells.Item(6,2).Value = $wsTotal.Cells.Item(5,2).Value2
$wsTotal.Cells.Item(6,3).Value = $wsTotal.Cells.Item(5,3).Value2
$wsTotal.Cells.Item(6,4).Value = $wsTotal.Cells.Item(5,4).Value2

$wsTotal.Cells.Item(5,2).Value = $wsTotal.Cells.Item(4,2).Value2
$wsTotal.Cells.Item(5,3).Value = $wsTotal.Cells.Item(4,3).Value2
$wsTotal.Cells.Item(5,4).Value = $wsTotal.Cells.Item(4,4).Value2

$wsTotal.Cells.Item(4,2).Value = $wsTotal.Cells.Item(3,2).Value2
$wsTotal.Cells.Item(4,3).Value = $wsTotal.Cells.Item(3,3).Value2
$wsTotal.Cells.Item(4,4).Value = $wsTotal.Cells.Item(3,4).Value2

$wsTotal.Cells.Item(3,2).Value = $wsTotal.Cells.Item(2,2).Value2
$wsTotal.Cells.Item(3,3).Value = $wsTotal.Cells.Item(2,3).Value2
$wsTotal.Cells.Item(3,4).Value = $wsTotal.Cells.Item(2,4).Value2

# New row 2: the 2022-Q1 summary.
$wsTotal.Cells.Item(2,2).Value = "2022-Q1"
$wsTotal.Cells.Item(2,3).Value = 5
$wsTotal.Cells.Item(2,4).Value = 0.78

# Renumber the zero-based index column for all six rows.
$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(6,1).Value = 4
